$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 514.41174
$ws.Range("I19").Value = 456.7857
$ws.Range("J19").Value = 783.3333
$ws.Range("K19").Value = 456.7857
$ws.Range("L19").Value = 783.3333
$ws.Range("M19").Value = -281.7857
$ws.Range("N19").Value = -1133.3333
$ws.Range("I51").Value = 2500
$ws.Range("J51").Value = 2800
$ws.Range("K51").Value = 2500
$ws.Range("L51").Value = 2800
$ws.Range("M51").Value = -2016
$ws.Range("N51").Value = -3768
$ws.Range("H129").Value = 1146.78
$ws.Range("I129").Value = 523.1429000000001
$ws.Range("J129").Value = 1248.3024
$ws.Range("K129").Value = 1569.4287
$ws.Range("L129").Value = 3744.9072
$ws.Range("M129").Value = 3430.5713
$ws.Range("N129").Value = -13744.9072
$ws.Range("H137").Value = 1109.9375
$ws.Range("I137").Value = 977.4838999999999
$ws.Range("J137").Value = 1351.4706
$ws.Range("K137").Value = 2932.4517
$ws.Range("L137").Value = 4054.4118
$ws.Range("M137").Value = -382.4516999999996
$ws.Range("N137").Value = -9154.4118

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H123").Value = 24288
$ws.Range("J123").Value = 24288
$ws.Range("L123").Value = 24288
$ws.Range("N123").Value = -34088

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 49589.24
$ws.Range("I20").Value = 78760.69500000001
$ws.Range("J20").Value = 2185.625
$ws.Range("K20").Value = 78760.69500000001
$ws.Range("L20").Value = 2185.625
$ws.Range("M20").Value = -78513.69500000001
$ws.Range("N20").Value = -2679.625
$ws.Range("H37").Value = 1668.2858
$ws.Range("I37").Value = 344.5
$ws.Range("J37").Value = 3433.3333
$ws.Range("K37").Value = 344.5
$ws.Range("L37").Value = 3433.3333
$ws.Range("M37").Value = -207.5
$ws.Range("N37").Value = -3707.3333
$ws.Range("H80").Value = 2645666.2
$ws.Range("I80").Value = 7407453
$ws.Range("J80").Value = 228.88889
$ws.Range("K80").Value = 7407453
$ws.Range("L80").Value = 228.88889
$ws.Range("M80").Value = -7406455
$ws.Range("N80").Value = -2224.88889
$ws.Range("H83").Value = 2645666.2
$ws.Range("I83").Value = 7407453
$ws.Range("J83").Value = 228.88889
$ws.Range("K83").Value = 37037265
$ws.Range("L83").Value = 1144.44445
$ws.Range("M83").Value = -37032273
$ws.Range("N83").Value = -11128.44445
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1934.2
$ws.Range("I31").Value = 1419.5946
$ws.Range("J31").Value = 4314.25
$ws.Range("K31").Value = 1419.5946
$ws.Range("L31").Value = 4314.25
$ws.Range("M31").Value = -1124.5946
$ws.Range("N31").Value = -4904.25
$ws.Range("H34").Value = 1934.2
$ws.Range("I34").Value = 1419.5946
$ws.Range("J34").Value = 4314.25
$ws.Range("K34").Value = 1419.5946
$ws.Range("L34").Value = 4314.25
$ws.Range("M34").Value = -1217.5946
$ws.Range("N34").Value = -4718.25
$ws.Range("H134").Value = 1482.5306
$ws.Range("I134").Value = 1310.6666
$ws.Range("K134").Value = 3931.9998
$ws.Range("M134").Value = -1396.9998

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2340.5715
$ws.Range("I132").Value = 1334.6666
$ws.Range("J132").Value = 2614.9092
$ws.Range("K132").Value = 12011.9994
$ws.Range("L132").Value = 23534.1828
$ws.Range("M132").Value = -9481.999400000001
$ws.Range("N132").Value = -28594.1828
$ws.Range("H133").Value = 4308.4443
$ws.Range("I133").Value = 1340.1111
$ws.Range("J133").Value = 7276.778
$ws.Range("K133").Value = 4020.3333
$ws.Range("L133").Value = 21830.334
$ws.Range("M133").Value = 1039.6667
$ws.Range("N133").Value = -31950.334
$ws.Range("H134").Value = 3856.0625
$ws.Range("I134").Value = 2004.8889
$ws.Range("J134").Value = 6236.143
$ws.Range("K134").Value = 6014.6667
$ws.Range("L134").Value = 18708.429
$ws.Range("M134").Value = -944.6666999999998
$ws.Range("N134").Value = -28848.429
$ws.Range("H138").Value = 2750.3333
$ws.Range("I138").Value = 1147.8
$ws.Range("K138").Value = 3443.4
$ws.Range("M138").Value = 1696.6

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 30000
$ws.Range("J62").Value = 30000
$ws.Range("L62").Value = 30000
$ws.Range("N62").Value = -31372
$ws.Range("H65").Value = 30000
$ws.Range("J65").Value = 30000
$ws.Range("L65").Value = 90000
$ws.Range("N65").Value = -96864
$ws.Range("H70").Value = 6541.136
$ws.Range("I70").Value = 5858
$ws.Range("J70").Value = 7360.9
$ws.Range("K70").Value = 5858
$ws.Range("L70").Value = 7360.9
$ws.Range("M70").Value = -5588
$ws.Range("N70").Value = -7900.9
$ws.Range("H73").Value = 6541.136
$ws.Range("I73").Value = 5858
$ws.Range("J73").Value = 7360.9
$ws.Range("K73").Value = 5858
$ws.Range("L73").Value = 7360.9
$ws.Range("M73").Value = -4922
$ws.Range("N73").Value = -9232.9
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 9456.552
$ws.Range("I61").Value = 11340.913
$ws.Range("J61").Value = 2233.1667
$ws.Range("K61").Value = 11340.913
$ws.Range("L61").Value = 2233.1667
$ws.Range("M61").Value = -11138.913
$ws.Range("N61").Value = -2637.1667
$ws.Range("H100").Value = 2197.1667
$ws.Range("I100").Value = 1601.5
$ws.Range("J100").Value = 2495
$ws.Range("K100").Value = 1601.5
$ws.Range("L100").Value = 2495
$ws.Range("M100").Value = -1060.5
$ws.Range("N100").Value = -3577
$ws.Range("H113").Value = 9456.552
$ws.Range("I113").Value = 11340.913
$ws.Range("J113").Value = 2233.1667
$ws.Range("K113").Value = 11340.913
$ws.Range("L113").Value = 2233.1667
$ws.Range("M113").Value = -9170.913
$ws.Range("N113").Value = -6573.1667
$ws.Range("H132").Value = 3500.6128
$ws.Range("I132").Value = 2955.2917
$ws.Range("K132").Value = 8865.875100000001
$ws.Range("M132").Value = -6335.875100000001

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 26924.834
$ws.Range("J123").Value = 26924.834
$ws.Range("L123").Value = 26924.834
$ws.Range("N123").Value = -36724.834
